$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.74"
$ws.Range("E2").Value = "'-0.64%"
$ws.Range("G2").Value = "'16"
$ws.Range("D3").Value = "'31.48"
$ws.Range("E3").Value = "'-2.69%"
$ws.Range("G3").Value = "'16"
$ws.Range("D4").Value = "'5.148"
$ws.Range("E4").Value = "'-2.83%"
$ws.Range("G4").Value = "'16"
$ws.Range("D5").Value = "'0.07399"
$ws.Range("E5").Value = "'-1.23%"
$ws.Range("G5").Value = "'16"
$ws.Range("D6").Value = "'1.829"
$ws.Range("E6").Value = "'22.91%"
$ws.Range("G6").Value = "'16"
$ws.Range("D7").Value = "'7.877"
$ws.Range("E7").Value = "'0.91%"
$ws.Range("G7").Value = "'16"
$ws.Range("D8").Value = "'3.770"
$ws.Range("E8").Value = "'-0.80%"
$ws.Range("G8").Value = "'16"
$ws.Range("D9").Value = "'0.9281"
$ws.Range("E9").Value = "'0.87%"
$ws.Range("G9").Value = "'16"
$ws.Range("D10").Value = "'0.1711"
$ws.Range("E10").Value = "'0.51%"
$ws.Range("G10").Value = "'16"
$ws.Range("D11").Value = "'0.07320"
$ws.Range("E11").Value = "'-5.97%"
$ws.Range("G11").Value = "'16"
$ws.Range("D12").Value = "'0.08138"
$ws.Range("E12").Value = "'1.10%"
$ws.Range("G12").Value = "'16"
$ws.Range("D13").Value = "'0.03041"
$ws.Range("E13").Value = "'0.46%"
$ws.Range("G13").Value = "'16"
$ws.Range("D14").Value = "'0.09944"
$ws.Range("E14").Value = "'0.57%"
$ws.Range("G14").Value = "'16"
$ws.Range("D15").Value = "'0.001498"
$ws.Range("E15").Value = "'0.34%"
$ws.Range("G15").Value = "'16"
$ws.Range("D16").Value = "'0.006067"
$ws.Range("E16").Value = "'-6.30%"
$ws.Range("G16").Value = "'16"
$ws.Range("D17").Value = "'3.471"
$ws.Range("E17").Value = "'-0.03%"
$ws.Range("G17").Value = "'16"
$ws.Range("D18").Value = "'2.223"
$ws.Range("E18").Value = "'-0.23%"
$ws.Range("G18").Value = "'16"
$ws.Range("D19").Value = "'0.3260"
$ws.Range("E19").Value = "'-2.02%"
$ws.Range("G19").Value = "'16"
$ws.Range("D20").Value = "'0.1347"
$ws.Range("E20").Value = "'0.27%"
$ws.Range("G20").Value = "'16"
$ws.Range("D21").Value = "'4.616"
$ws.Range("E21").Value = "'2.98%"
$ws.Range("G21").Value = "'16"
$ws.Range("D22").Value = "'0.04656"
$ws.Range("E22").Value = "'1.19%"
$ws.Range("G22").Value = "'16"
$ws.Range("D23").Value = "'0.1578"
$ws.Range("E23").Value = "'-2.59%"
$ws.Range("G23").Value = "'16"
$ws.Range("D24").Value = "'0.001217"
$ws.Range("E24").Value = "'0.09%"
$ws.Range("G24").Value = "'16"
$ws.Range("D25").Value = "'0.004496"
$ws.Range("E25").Value = "'1.50%"
$ws.Range("G25").Value = "'16"
$ws.Range("E26").Value = "'-7.41%"
$ws.Range("G26").Value = "'16"
$ws.Range("E27").Value = "'7.61%"
$ws.Range("G27").Value = "'16"
$ws.Range("G28").Value = "'16"
$ws.Range("G29").Value = "'16"
$ws.Range("G30").Value = "'16"
$ws.Range("G31").Value = "'16"
$ws.Range("G32").Value = "'16"
$ws.Range("G33").Value = "'16"
$ws.Range("G34").Value = "'16"
$ws.Range("G35").Value = "'16"
$ws.Range("G36").Value = "'16"
$ws.Range("G37").Value = "'16"
$ws.Range("G38").Value = "'16"
$ws.Range("D39").Value = "'0.01724"
$ws.Range("E39").Value = "'-3.49%"
$ws.Range("G39").Value = "'16"
$ws.Range("D40").Value = "'0.04517"
$ws.Range("E40").Value = "'-0.75%"
$ws.Range("G40").Value = "'16"
$ws.Range("D41").Value = "'0.007118"
$ws.Range("E41").Value = "'-1.57%"
$ws.Range("G41").Value = "'16"
$ws.Range("D42").Value = "'0.1351"
$ws.Range("E42").Value = "'0.77%"
$ws.Range("G42").Value = "'16"
$ws.Range("D43").Value = "'0.002144"
$ws.Range("E43").Value = "'-2.12%"
$ws.Range("G43").Value = "'16"
$ws.Range("D44").Value = "'0.01086"
$ws.Range("E44").Value = "'-14.28%"
$ws.Range("G44").Value = "'16"
$ws.Range("D45").Value = "'0.00006225"
$ws.Range("E45").Value = "'3.06%"
$ws.Range("G45").Value = "'16"
$ws.Range("E46").Value = "'-21.35%"
$ws.Range("G46").Value = "'16"
$ws.Range("E47").Value = "'159.96%"
$ws.Range("G47").Value = "'16"
$ws.Range("G48").Value = "'16"
$ws.Range("G49").Value = "'16"
$ws.Range("G50").Value = "'16"
$ws.Range("G51").Value = "'16"
